$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (11) and populate its header (row 3)
# with the new "Hyperscaler" value. This pushes the existing K:P columns
# (and their widths / merged ranges / shared-string references) one slot
# to the right, to L:Q.
$ws.Columns("K").Insert() | Out-Null
$ws.Range("K3").Value = "Hyperscaler"

# Restore the explicit column width for the freshly inserted column K
# (closest value reachable through the ColumnWidth property's pixel
# rounding to the original authored width of 10.4609375 characters).
$ws.Columns("K").ColumnWidth = 9.666666666666666

# Keep the header-row selection in sync with the widened merged range
# (A1:P1 -> A1:Q1), matching how Excel updates the active selection
# after the new column pushes the merge area out by one column.
$ws.Range("A1:Q1").Select() | Out-Null
